# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates the "K" column (G) values for rows 2-24 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 1
    6  = 2
    7  = 0
    8  = 1
    9  = 0
    10 = 2
    11 = 2
    12 = 1
    13 = 1
    14 = 1
    15 = 2
    16 = 1
    17 = 1
    18 = 2
    19 = 2
    20 = 2
    21 = 1
    22 = 2
    23 = 0
    24 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
